# INFO2_MOD07-act7.docx — "Update activity nextTag error"
#
# Fixes two content errors in the activity instructions:
#   1. The paragraph describing how the field is modified mistakenly
#      referred to it as "nextTag" instead of its real name, "tag"
#      (the field is introduced/used as "tag" everywhere else in the
#      document, e.g. "changeTag", "el valor de la variable tag").
#   2. The sample-output table for AndroidDemo2 listed non-prime
#      numbers (25 and 35) where the program's changeTag()/isPrime()
#      logic only ever produces primes; the expected sample output is
#      corrected to the next primes after 25 and 31 respectively
#      (29, 31, 37, 41).

$d = $word.ActiveDocument

# --- Fix 1: "...el valor de la variable nextTag será modificado..." -> "tag" ---
$range1 = $d.Content
$found1 = $range1.Find.Execute("nextTag", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "tag", 2)
if (-not $found1) {
    throw "Could not find 'nextTag' to fix."
}

# --- Fix 2: correct the sample output table (2nd table) to real primes ---
$table = $d.Tables.Item(2)
$cell = $table.Cell(1, 1)

$pairs = @(
    @("Jack25", "Jack29"),
    @("Bob29",  "Bob31"),
    @("Jack31", "Jack37"),
    @("Bob35",  "Bob41")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $cellRange = $cell.Range
    $found = $cellRange.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Could not find '$old' in the sample output table."
    }
}
